$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (19) with the new problem entry, mirroring the existing
# Problem / Class column layout.
$ws.Range("A19").Value = "Find first non-repeating element in a given Array of integers"
$ws.Range("B19").Value = "FirstNonRepeatingElement"

# Match the author's final selection after adding the row.
$ws.Range("B20").Select()
